# The recorded diff only reorders the xmlns:* attribute declarations on the
# six <ma14:wrappingTextBoxFlag .../> elements living inside each shape's
# <a:extLst><a:ext uri="{C572A759-...}"> (Mac PowerPoint's legacy "text box
# wraps text" flag, on "PPT统一标准" in slides 1-3, "文本框 6" in
# slideLayout15, and the Title/Body placeholders of slideMaster1). The
# attribute *set* and the val="1" payload are byte-identical before/after -
# only the order of the (redundant) namespace declarations changes, which is
# not semantically meaningful in XML. It is not tied to any PowerPoint
# object-model property (word-wrap is modelled through bodyPr/@wrap, a
# separate/modern mechanism) so there is nothing to toggle here without
# fabricating an unrelated, visible side effect. Touch nothing so the rest
# of the deck stays pixel-for-pixel identical to the source file.
$p = $ppt.ActivePresentation
